$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) retains text formatting so numeric-looking
# strings like "1.00" or "72.080.75" are not coerced to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "72.080.75"
$ws.Range("E2").Value = "  +3.87%  "

$ws.Range("D3").Value = "3.661.79"
$ws.Range("E3").Value = "  +7.39%  "

$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").Value = "593.97"
$ws.Range("E5").Value = "  +1.07%  "

$ws.Range("D6").Value = "181.90"
$ws.Range("E6").Value = "  +0.41%  "

$ws.Range("D7").Value = "3.654.77"
$ws.Range("E7").Value = "  +7.42%  "

$ws.Range("E8").Value = "  +1.77%  "

$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.03%  "

$ws.Range("E10").Value = "  +3.77%  "

$ws.Range("D11").Value = "0.608"
$ws.Range("E11").Value = "  +2.44%  "

$ws.Range("D12").Value = "50.13"
$ws.Range("E12").Value = "  +3.06%  "

$ws.Range("E13").Value = "  +1.74%  "

$ws.Range("D14").Value = "694.59"
$ws.Range("E14").Value = "  +1.40%  "

$ws.Range("D15").Value = "4.245.97"
$ws.Range("E15").Value = "  +7.25%  "

$ws.Range("D16").Value = "9.03"
$ws.Range("E16").Value = "  +4.25%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.699.29"
$ws.Range("E17").Value = "  +8.64%  "

$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "72.168.81"
$ws.Range("E18").Value = "  +3.81%  "

$ws.Range("E19").Value = "  +2.40%  "

$ws.Range("D20").Value = "18.61"
$ws.Range("E20").Value = "  +4.58%  "

$ws.Range("E21").Value = "  +2.47%  "

$ws.Range("D22").Value = "0.939"
$ws.Range("E22").Value = "  +3.15%  "

$ws.Range("D23").Value = "5.84"
$ws.Range("E23").Value = "  +8.13%  "

$ws.Range("D24").Value = "17.96"
$ws.Range("E24").Value = "  +4.18%  "

$ws.Range("D25").Value = "104.01"
$ws.Range("E25").Value = "  +0.78%  "

$ws.Range("E26").Value = "  +2.26%  "

$ws.Range("E27").Value = "  +5.04%  "

$ws.Range("D28").Value = "10.04"
$ws.Range("E28").Value = "  +3.14%  "

$ws.Range("D29").Value = "35.37"
$ws.Range("E29").Value = "  +3.55%  "

$ws.Range("E30").Value = "  +3.59%  "

$ws.Range("E31").Value = "  +5.67%  "

$ws.Range("D32").Value = "4.18"
$ws.Range("E32").Value = "  +16.28%  "

$ws.Range("D33").Value = "583.88"
$ws.Range("E33").Value = "  +3.59%  "

$ws.Range("D34").Value = "11.39"
$ws.Range("E34").Value = "  +1.82%  "

$ws.Range("E35").Value = "  +2.95%  "

$ws.Range("D36").Value = "59.69"
$ws.Range("E36").Value = "  +2.27%  "

$ws.Range("E37").Value = "  +0.02%  "

$ws.Range("D38").Value = "3.681.04"
$ws.Range("E38").Value = "  +0.21%  "

$ws.Range("E39").Value = "  +1.87%  "

$ws.Range("E40").Value = "  +0.00%  "

$ws.Range("D41").Value = "0.0₃0775"
$ws.Range("E41").Value = "  +7.29%  "

$ws.Range("D42").Value = "3.43"
$ws.Range("E42").Value = "  +4.75%  "

$ws.Range("E43").Value = "  +8.70%  "

$ws.Range("E44").Value = "  +4.01%  "

$ws.Range("E45").Value = "  +3.05%  "

$ws.Range("E46").Value = "  +3.92%  "

$ws.Range("E47").Value = "  +6.36%  "

$ws.Range("E48").Value = "  +2.26%  "

$ws.Range("E49").Value = "  +3.66%  "

$ws.Range("E50").Value = "  -0.08%  "

$ws.Range("E51").Value = "  +14.52%  "

# Restore the original (default) cell style now that values are set,
# so only the number formatting used for text-coercion is removed.
$ws.Range("D2:D51").Style = "Normal"
